# Auto-generated edit script applying the diff changes to Sargatanas_Profits workbook
# Updates market-data columns (H-N) for specific leve rows across multiple sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 42055970
$ws.Range("I70").Value = 17339400
$ws.Range("J70").Value = 75760380
$ws.Range("K70").Value = 52018200
$ws.Range("L70").Value = 227281140
$ws.Range("M70").Value = -52017930
$ws.Range("N70").Value = -227281680

$ws.Range("H73").Value = 42055970
$ws.Range("I73").Value = 17339400
$ws.Range("J73").Value = 75760380
$ws.Range("K73").Value = 52018200
$ws.Range("L73").Value = 227281140
$ws.Range("M73").Value = -52017264
$ws.Range("N73").Value = -227283012

$ws.Range("H107").Value = 44232480
$ws.Range("I107").Value = 28126580
$ws.Range("J107").Value = 70001930
$ws.Range("K107").Value = 28126580
$ws.Range("L107").Value = 70001930
$ws.Range("M107").Value = -28124660
$ws.Range("N107").Value = -70005770

$ws.Range("H113").Value = 115391780
$ws.Range("I113").Value = 2782
$ws.Range("K113").Value = 2782
$ws.Range("M113").Value = 472

$ws.Range("H132").Value = 1684.5476
$ws.Range("I132").Value = 1635.7428
$ws.Range("J132").Value = 1928.5714
$ws.Range("K132").Value = 4907.2284
$ws.Range("L132").Value = 5785.7142
$ws.Range("M132").Value = -2377.2284
$ws.Range("N132").Value = -10845.7142

$ws.Range("H137").Value = 3865.7917
$ws.Range("J137").Value = 3410.4707
$ws.Range("L137").Value = 10231.4121
$ws.Range("N137").Value = -15331.4121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1871831.6
$ws.Range("I32").Value = 2122846
$ws.Range("K32").Value = 2122846
$ws.Range("M32").Value = -2122559

$ws.Range("H45").Value = 3398.818
$ws.Range("I45").Value = 2474
$ws.Range("J45").Value = 4039.077
$ws.Range("K45").Value = 2474
$ws.Range("L45").Value = 4039.077
$ws.Range("M45").Value = -2097
$ws.Range("N45").Value = -4793.077

$ws.Range("H61").Value = 8997.684999999999
$ws.Range("I61").Value = 1994
$ws.Range("K61").Value = 1994
$ws.Range("M61").Value = -1782

$ws.Range("H97").Value = 3969005.8
$ws.Range("I97").Value = 789.15
$ws.Range("J97").Value = 83333336
$ws.Range("K97").Value = 789.15
$ws.Range("L97").Value = 83333336
$ws.Range("M97").Value = -293.15
$ws.Range("N97").Value = -83334328

$ws.Range("H110").Value = 41668776
$ws.Range("I110").Value = 1772.2
$ws.Range("J110").Value = 111113780
$ws.Range("K110").Value = 1772.2
$ws.Range("L110").Value = 111113780
$ws.Range("M110").Value = 272.8
$ws.Range("N110").Value = -111117870

$ws.Range("H132").Value = 4207
$ws.Range("I132").Value = 2038.3778
$ws.Range("K132").Value = 6115.1334
$ws.Range("M132").Value = -3585.1334

$ws.Range("H136").Value = 8997.684999999999
$ws.Range("I136").Value = 1994
$ws.Range("K136").Value = 5982
$ws.Range("M136").Value = -3432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6412305.5
$ws.Range("I20").Value = 8773925
$ws.Range("K20").Value = 8773925
$ws.Range("M20").Value = -8773678

$ws.Range("H26").Value = 31827.143
$ws.Range("I26").Value = 22981.2
$ws.Range("K26").Value = 22981.2
$ws.Range("M26").Value = -22689.2

$ws.Range("H28").Value = 53841.5
$ws.Range("J28").Value = 53841.5
$ws.Range("L28").Value = 53841.5
$ws.Range("N28").Value = -54429.5

$ws.Range("H86").Value = 48117230
$ws.Range("I86").Value = 10914662
$ws.Range("J86").Value = 333336930
$ws.Range("K86").Value = 10914662
$ws.Range("L86").Value = 333336930
$ws.Range("M86").Value = -10913539
$ws.Range("N86").Value = -333339176

$ws.Range("H89").Value = 48117230
$ws.Range("I89").Value = 10914662
$ws.Range("J89").Value = 333336930
$ws.Range("K89").Value = 54573310
$ws.Range("L89").Value = 1666684650
$ws.Range("M89").Value = -54567694
$ws.Range("N89").Value = -1666695882

$ws.Range("H94").Value = 1093.2632
$ws.Range("J94").Value = 1809.1111
$ws.Range("L94").Value = 1809.1111
$ws.Range("N94").Value = -2711.1111

$ws.Range("H96").Value = 20553.857

$ws.Range("H134").Value = 5259.75
$ws.Range("I134").Value = 1933.2727
$ws.Range("K134").Value = 5799.8181
$ws.Range("M134").Value = -3264.8181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5854071
$ws.Range("I31").Value = 2514
$ws.Range("J31").Value = 13899962
$ws.Range("K31").Value = 2514
$ws.Range("L31").Value = 13899962
$ws.Range("M31").Value = -2219
$ws.Range("N31").Value = -13900552

$ws.Range("H34").Value = 5854071
$ws.Range("I34").Value = 2514
$ws.Range("J34").Value = 13899962
$ws.Range("K34").Value = 2514
$ws.Range("L34").Value = 13899962
$ws.Range("M34").Value = -2312
$ws.Range("N34").Value = -13900366

$ws.Range("H58").Value = 6761167.5
$ws.Range("I58").Value = 10640199
$ws.Range("K58").Value = 10640199
$ws.Range("M58").Value = -10639996

$ws.Range("H132").Value = 5337527
$ws.Range("I132").Value = 2210.5715
$ws.Range("K132").Value = 6631.7145
$ws.Range("M132").Value = -4101.7145

$ws.Range("H134").Value = 4868.162
$ws.Range("I134").Value = 2661.9524
$ws.Range("K134").Value = 7985.8572
$ws.Range("M134").Value = -5450.8572

$ws.Range("H136").Value = 6761167.5
$ws.Range("I136").Value = 10640199
$ws.Range("K136").Value = 31920597
$ws.Range("M136").Value = -31918047

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7383.1665
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 7383.1665
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 22149.4995
$ws.Range("N5").Value = -22373.4995
$ws.Range("M5").ClearContents()

$ws.Range("H50").Value = 66666828
$ws.Range("I50").Value = 83333530
$ws.Range("J50").Value = 20
$ws.Range("K50").Value = 250000590
$ws.Range("L50").Value = 60
$ws.Range("M50").Value = -250000109
$ws.Range("N50").Value = -1022

$ws.Range("H53").Value = 66666828
$ws.Range("I53").Value = 83333530
$ws.Range("J53").Value = 20
$ws.Range("K53").Value = 250000590
$ws.Range("L53").Value = 60
$ws.Range("M53").Value = -250000109
$ws.Range("N53").Value = -1022

$ws.Range("H55").Value = 71673880
$ws.Range("J55").Value = 5890752.5
$ws.Range("L55").Value = 17672257.5
$ws.Range("N55").Value = -17672611.5

$ws.Range("H111").Value = 11114.833
$ws.Range("I111").Value = 11114.833
$ws.Range("K111").Value = 33344.499
$ws.Range("M111").Value = -30277.499

$ws.Range("H131").Value = 1954.641
$ws.Range("J131").Value = 2380.4
$ws.Range("L131").Value = 7141.200000000001
$ws.Range("N131").Value = -17221.2

$ws.Range("H135").Value = 7383.1665
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 7383.1665
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 66448.4985
$ws.Range("N135").Value = -71518.4985
$ws.Range("M135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 58833972
$ws.Range("I70").Value = 166675250
$ws.Range("K70").Value = 166675250
$ws.Range("M70").Value = -166674980

$ws.Range("H73").Value = 58833972
$ws.Range("I73").Value = 166675250
$ws.Range("K73").Value = 166675250
$ws.Range("M73").Value = -166674314

$ws.Range("H113").Value = 5835.5107
$ws.Range("I113").Value = 2570.8635
$ws.Range("J113").Value = 8708.4
$ws.Range("K113").Value = 2570.8635
$ws.Range("L113").Value = 8708.4
$ws.Range("M113").Value = -400.8634999999999
$ws.Range("N113").Value = -13048.4

$ws.Range("H122").Value = 4216673
$ws.Range("I122").Value = 8955473
$ws.Range("K122").Value = 26866419
$ws.Range("M122").Value = -26863969

$ws.Range("H132").Value = 6523.3794
$ws.Range("I132").Value = 2636.25
$ws.Range("K132").Value = 7908.75
$ws.Range("M132").Value = -5378.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4385
$ws.Range("I40").Value = 2111.875
$ws.Range("K40").Value = 2111.875
$ws.Range("M40").Value = -1975.875

$ws.Range("H122").Value = 6607
$ws.Range("I122").Value = 3408.5
$ws.Range("J122").Value = 8083.231
$ws.Range("K122").Value = 10225.5
$ws.Range("L122").Value = 24249.693
$ws.Range("M122").Value = -7775.5
$ws.Range("N122").Value = -29149.693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 146414.61
$ws.Range("I122").Value = 269506.8
$ws.Range("J122").Value = 4385.154
$ws.Range("K122").Value = 808520.3999999999
$ws.Range("L122").Value = 13155.462
$ws.Range("M122").Value = -806070.3999999999
$ws.Range("N122").Value = -18055.462

$ws.Range("H126").Value = 3126.8333
$ws.Range("I126").Value = 1464
$ws.Range("K126").Value = 4392
$ws.Range("M126").Value = -1922

$ws.Range("H132").Value = 12204692
$ws.Range("I132").Value = 14709264
$ws.Range("K132").Value = 44127792
$ws.Range("M132").Value = -44125262

$ws.Range("H136").Value = 26346740
$ws.Range("I136").Value = 52632560
$ws.Range("K136").Value = 157897680
$ws.Range("M136").Value = -157895130
